$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first column (old "runoff" index column with values 0..3) and
# shift everything left. After this, the old "peak_runoff_rate" column
# (old B, containing the runoff data) becomes column A, the old C/D/E
# columns become B/C/D, and the "Destore-Perv" column (old F) becomes
# column E.
$ws.Columns.Item(1).Delete() | Out-Null

# Consolidate the remaining data: keep the runoff header/values in column A
# and move the "Destore-Perv" header/values (now in column E) into column B,
# right next to "runoff".
$ws.Range("B1").Value = $ws.Range("E1").Value()
$ws.Range("B2:B5").Value = $ws.Range("E2:E5").Value()

# Drop the now-unused peak_runoff_rate / infiltration / evaporation columns
# (and the now-duplicated Destore-Perv column). Their header cells keep
# their formatting, but all content is cleared out.
$ws.Range("C1:E5").ClearContents() | Out-Null

# Leave the same cell selected as in the authored workbook.
$ws.Range("E4").Select() | Out-Null
